# "Miglioramenti e pulizia generale"
#
# 1) Cell A3's text changes from the literal "${i}" to "${String}".
# 2) A new comment/note is added on A3 (mirroring the existing A1 note),
#    documenting the jx:each() template directives that replaced the old
#    jx:area() directive.
# 3) The worksheet's active selection moves from A4 to G8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Update the templated placeholder text in A3 ---------------------
# NOTE: use single-quoted strings so PowerShell does not try to expand
# "${i}" / "${String}" as variable references.
$ws.Range("A3").Value = '${String}'

# --- 2) Document the new jx:each() directives via a cell comment --------
$line1 = 'jx:each(items="master" var="items" lastCell="A3" direction="RIGHT")'
$line2 = 'jx:each(items="items" var="String" lastCell="A3" direction="DOWN")'
$commentText = "Author:" + [char]10 + $line1 + [char]10 + $line2

$comment = $ws.Range("A3").AddComment($commentText)
$comment.Visible = $false

# --- 3) Move the active selection from A4 to G8 --------------------------
$ws.Range("G8").Select()
